$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - copy style from an existing header cell (E1) so F1 gets
# the same bold / bordered / centered formatting.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Data cells - plain inline strings, no special style.
$ws.Range("F2").Value = "2021-10-05 13:41:06.800827"
$ws.Range("F3").Value = "2021-10-05 13:41:06.800837"
$ws.Range("F4").Value = "2021-10-05 13:41:06.800841"
$ws.Range("F5").Value = "2021-10-05 13:41:06.800844"
